$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 5 for columns D, J, K, L, M, P
$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow5 = $ws.Range($col + "5")

    $val2 = $cellRow2.Value2
    $val5 = $cellRow5.Value2

    $cellRow2.Value2 = $val5
    $cellRow5.Value2 = $val2
}
